$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.136.21'
$ws.Range("E2").Value = '  -0.31%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.842.44'
$ws.Range("E3").Value = '  -0.41%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.73'
$ws.Range("E5").Value = '  -1.72%  '

$ws.Range("E6").Value = '  -1.80%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.0000'
$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  -1.41%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07470'
$ws.Range("E9").Value = '  -3.29%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.14'
$ws.Range("E10").Value = '  -1.63%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07669'
$ws.Range("E11").Value = '  -1.91%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.845.20'
$ws.Range("E12").Value = '  -0.18%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.064'
$ws.Range("E13").Value = '  -1.24%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6840'
$ws.Range("E14").Value = '  -0.36%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '87.57'
$ws.Range("E15").Value = '  -5.86%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.176'
$ws.Range("E16").Value = '  -6.83%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.129.03'
$ws.Range("E17").Value = '  -0.27%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008170'
$ws.Range("E18").Value = '  -1.75%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '2.081.46'
$ws.Range("E19").Value = '  -0.64%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '228.26'
$ws.Range("E20").Value = '  -5.57%  '

$ws.Range("E21").Value = '  -1.50%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9994'
$ws.Range("E22").Value = '  -0.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.411'
$ws.Range("E23").Value = '  -1.51%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9999'
$ws.Range("E24").Value = '  +0.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1455'
$ws.Range("E25").Value = '  -3.58%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '160.08'
$ws.Range("E26").Value = '  +0.50%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.765'
$ws.Range("E27").Value = '  -0.72%  '

$ws.Range("E28").Value = '  -1.07%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.513'
$ws.Range("E29").Value = '  -2.09%  '

$ws.Range("E30").Value = '  +1.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.140'
$ws.Range("E31").Value = '  -1.00%  '

$ws.Range("E32").Value = '  -0.19%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05239'
$ws.Range("E33").Value = '  +2.30%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7646'
$ws.Range("E34").Value = '  -3.68%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.850'
$ws.Range("E35").Value = '  -1.68%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.136'
$ws.Range("E36").Value = '  -1.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.686'
$ws.Range("E37").Value = '  -0.27%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.312.71'
$ws.Range("E38").Value = '  -0.75%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01840'
$ws.Range("E39").Value = '  -1.90%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.726'
$ws.Range("E40").Value = '  +0.39%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9347'
$ws.Range("E41").Value = '  -1.48%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '105.19'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.795'
$ws.Range("E43").Value = '  -3.59%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9994'
$ws.Range("E44").Value = '  -0.09%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").Value = '1.982.59'
$ws.Range("E45").Value = '  -0.36%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '0.5201'
$ws.Range("E46").Value = '  +0.35%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '64.87'
$ws.Range("E47").Value = '  +1.11%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '9.547'
$ws.Range("E48").Value = '  -1.85%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '1.775'
$ws.Range("E49").Value = '  +0.48%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("B50").Value = 'XinFinNetwork'
$ws.Range("C50").Value = 'https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc'
$ws.Range("D50").Value = '0.07494'
$ws.Range("E50").Value = '  +18.82%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("B51").Value = 'Cronos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D51").Value = '0.05944'
$ws.Range("E51").Value = '  +0.86%  '
